# Append-edit for the "ランサーズ" (Lancers) sheet:
#  - two brand-new job postings are inserted right after the existing
#    header/data block at rows 10-11 (pushing the previous rows 10-13
#    down to rows 12-15)
#  - one more brand-new job posting is appended at the new last row (16)
#  - every data row's "取得日時" (fetched-at) timestamp in column A is
#    refreshed to the new scrape time

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-23 01:24:52"

# ---------------------------------------------------------------------
# 1) Make room for the two new postings: insert 2 blank rows at row 10,
#    pushing the old rows 10-13 down to rows 12-15 (data + formatting
#    shift together).
# ---------------------------------------------------------------------
$ws.Rows("10:11").Insert()

# ---------------------------------------------------------------------
# 2) Fill the two newly-inserted rows (10 and 11) with the new postings.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = $newTimestamp
$ws.Range("B10").Value = '【急募】n8nを使った請求書自動化プロジェクトの依頼'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5459128'
$ws.Range("G10").Value = 88
$ws.Range("H10").Value = '◆自動化'

$ws.Range("A11").Value = $newTimestamp
$ws.Range("B11").Value = '【急募】女性顧客向けチャットボット開発のプロを探しています!'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5458992'
$ws.Range("G11").Value = 68
$ws.Range("H11").Value = '◆開発'

# ---------------------------------------------------------------------
# 3) Append a brand-new row 16 with one more posting (note: it has no
#    "スキル概要" / column H value, same as the pre-existing Delphi row).
# ---------------------------------------------------------------------
$ws.Range("A16").Value = $newTimestamp
$ws.Range("B16").Value = '【電卓設計】ハードウェアとソフトウェアの専門家を募集!'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5459232'
$ws.Range("G16").Value = 18

# ---------------------------------------------------------------------
# 4) Refresh the "取得日時" timestamp for every other existing data row
#    (rows 2-9 kept their position; rows 12-15 are the shifted old
#    rows 10-13) so the whole sheet reflects the latest scrape run.
# ---------------------------------------------------------------------
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
foreach ($r in 12..15) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# ---------------------------------------------------------------------
# 5) Rebuild the URL hyperlinks for column F top to bottom. Row-insert
#    does not renumber the underlying relationship ids, so the cleanest
#    way to keep everything consistent is to clear every hyperlink on
#    the sheet and re-add them in row order. (Reading back `.Value`
#    through this host is unreliable, so use `.Value2` which returns
#    the real stored string.)
# ---------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value2
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}
